# Weekly update: Fruta / hortaliza, semanal
# Insert 6 new rows at row 598 (new price week 2022-05-25 / serial 44706)
# for "Comercializadora del Agro de Limari" - Tomate, pushing the existing
# rows 598-624 down to 604-630.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert six blank rows above the current row 598.
$ws.Range("A598:A603").EntireRow.Insert()

# Common (constant) values shared by every data row in this block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$categoriaId = 100112020
$categoria = "Tomate"
$unidad    = "$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"

$fecha = 44706

$newRows = @(
    @{ Row = 598; Variedad = "Larga vida"; Calidad = "Primera"; Volumen = 1100; PMin = 14000; PMax = 15000; PProm = 14500; PKg = 806 },
    @{ Row = 599; Variedad = "Larga vida"; Calidad = "Segunda"; Volumen = 800;  PMin = 12000; PMax = 13000; PProm = 12500; PKg = 694 },
    @{ Row = 600; Variedad = "Larga vida"; Calidad = "Tercera"; Volumen = 700;  PMin = 10000; PMax = 11000; PProm = 10500; PKg = 583 },
    @{ Row = 601; Variedad = "Semiduro";   Calidad = "Primera"; Volumen = 400;  PMin = 11000; PMax = 12000; PProm = 11500; PKg = 639 },
    @{ Row = 602; Variedad = "Semiduro";   Calidad = "Segunda"; Volumen = 400;  PMin = 9000;  PMax = 10000; PProm = 9500;  PKg = 528 },
    @{ Row = 603; Variedad = "Semiduro";   Calidad = "Tercera"; Volumen = 300;  PMin = 7000;  PMax = 8000;  PProm = 7500;  PKg = 417 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $r.Variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
